# ---------------------------------------------------------------------------
# "Finished button masher project" edit
#
#   1. Split "Print DocPacs, Sept. " into 3 runs, flagging "DocPacs" with a
#      spell-check proofErr pair.
#   2. Added new list items ("[J] Pseudocode Practice", "[J] Button Masher
#      Game", "Desktop Mount Installation") to the "Included Documentation"
#      and "Required Documentation" table cells.
#   3. Wrapped "formbar" in a spell-check proofErr pair.
#   4. Split "Serious contribution to a issue, ..." into 3 runs, flagging the
#      lone "a" with a spell-check proofErr pair.
#   5. Fixed a typo: "they way" -> "the way".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$WNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Clean-Text($s) {
    # Range.Text can carry a trailing paragraph mark (13) and/or a
    # end-of-cell marker (7) that plain .Trim() will not strip.
    return $s.TrimEnd([char]7, [char]13, [char]10).Trim()
}

function Find-ParaIndex($doc, $text, $startAt) {
    $count = $doc.Paragraphs.Count
    for ($i = $startAt; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ((Clean-Text $p.Range.Text) -eq $text) {
            return $i
        }
    }
    return -1
}

function New-ListParaXml($numId, $text) {
    return "<w:p $WNS><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"$numId`"/></w:numPr><w:ind w:right=`"240`"/><w:rPr><w:rFonts w:ascii=`"Abadi`" w:hAnsi=`"Abadi`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Abadi`" w:hAnsi=`"Abadi`"/></w:rPr><w:t>$text</w:t></w:r></w:p>"
}

# ---------------------------------------------------------------------------
# 1. "Print DocPacs, Sept. 29th" -- flag "DocPacs" as a spelling exception.
# ---------------------------------------------------------------------------
$printIdx = Find-ParaIndex $d "Print DocPacs, Sept. 29th" 1
if ($printIdx -gt 0) {
    $xml = '<w:p ' + $WNS + '>' +
        '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>' +
        '<w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr>' +
        '<w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/></w:rPr><w:t xml:space="preserve">Print </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/></w:rPr><w:t>DocPacs</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/></w:rPr><w:t xml:space="preserve">, Sept. </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/></w:rPr><w:t>29</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>' +
        '</w:p>'
    $d.Paragraphs($printIdx).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2a. "Included Documentation" cell (numId 39): the old "Desktop Mount
#     Installation" item becomes "[J] Pseudocode Practice", followed by new
#     "[J] Button Masher Game" and "Desktop Mount Installation" items, right
#     before the existing "Reflection" item.
# ---------------------------------------------------------------------------
$cell1Idx = Find-ParaIndex $d "Desktop Mount Installation" 1
if ($cell1Idx -gt 0) {
    $xml1 = (New-ListParaXml 39 "[J] Pseudocode Practice") +
            (New-ListParaXml 39 "[J] Button Masher Game") +
            (New-ListParaXml 39 "Desktop Mount Installation")
    $d.Paragraphs($cell1Idx).Range.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# 2b. "Required Documentation" cell (numId 37): add "Desktop Mount
#     Installation", "[J] Pseudocode Practice" and "[J] Button Masher Game"
#     right before the existing "Reflection" item.
# ---------------------------------------------------------------------------
$reqHeadingIdx = Find-ParaIndex $d "Required Documentation:" 1
$cell2ReflectionIdx = Find-ParaIndex $d "Reflection" $reqHeadingIdx
if ($cell2ReflectionIdx -gt 0) {
    $insertionRange = $d.Paragraphs($cell2ReflectionIdx).Range
    $insertionPoint = $d.Range($insertionRange.Start, $insertionRange.Start)
    $xml2 = (New-ListParaXml 37 "Desktop Mount Installation") +
            (New-ListParaXml 37 "[J] Pseudocode Practice") +
            (New-ListParaXml 37 "[J] Button Masher Game")
    $insertionPoint.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# 3. "You must contribute for the csmith1188/formbar and document your
#    contribution. ..." -- flag "formbar" as a spelling exception.
# ---------------------------------------------------------------------------
$formbarIdx = Find-ParaIndex $d "You must contribute for the csmith1188/formbar and document your contribution. A contribution is considered to be:" 1
if ($formbarIdx -gt 0) {
    $xml = '<w:p ' + $WNS + '>' +
        '<w:r><w:t xml:space="preserve">You must contribute for the csmith1188/</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>formbar</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> and document </w:t></w:r>' +
        '<w:r><w:t>your</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> cont</w:t></w:r>' +
        '<w:r><w:t>ribution. A contribution is considered to be:</w:t></w:r>' +
        '</w:p>'
    $d.Paragraphs($formbarIdx).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 4. "Serious contribution to a issue, discussions, or documentation" --
#    flag the lone "a" as a spelling exception.
# ---------------------------------------------------------------------------
$seriousIdx = Find-ParaIndex $d "Serious contribution to a issue, discussions, or documentation" 1
if ($seriousIdx -gt 0) {
    $xml = '<w:p ' + $WNS + '>' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="41"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Serious contribution to </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>a</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> issue, discussions, or documentation</w:t></w:r>' +
        '</w:p>'
    $d.Paragraphs($seriousIdx).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 5. Typo fix: "they way" -> "the way".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "How do you think they way you treat your workspace affects others in the school?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "How do you think the way you treat your workspace affects others in the school?",
    2) | Out-Null

Write-Output "edit complete"
